$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 8 into a new row 9 (this preserves cell formatting/style
# and the shared-string references for the text columns A-D).
$ws.Range("A8:F8").Copy($ws.Range("A9:F9"))

# Row 8 was missing the "danger_threshold" (F) check being lower than the
# "warning_threshold" (E) value; add the missing error row and fix the
# existing row's values: swap E8/F8 so E8=16.0 and F8=15.0.
$ws.Range("E8").Value2 = 16.0
$ws.Range("F8").Value2 = 15.0

# New row 9 keeps the same A-D values as row 8 (copied above) with
# E9=16.0 and F9=15.0.
$ws.Range("E9").Value2 = 16.0
$ws.Range("F9").Value2 = 15.0
